# The IAM032 / IAM033 test case rows were moved out of the "Test Cases"
# sheet (to ENWIAM), so delete those two rows entirely.
#
# Row 33 = IAM032 (TCID in column A), Row 34 = IAM033.
# Deleting row 33 twice removes both (row 34 shifts up to become the new
# row 33 after the first delete).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

$ws.Rows.Item(33).Delete()
$ws.Rows.Item(33).Delete()
